# Fruta / hortaliza, semanal
#
# Insert three new weekly price rows (grower "Patterson", week of 2023-01-05,
# serial 44931) for "Femacal de La Calera" / "Damasco" just above the
# existing row 203, pushing the old rows 203-212 down to 206-215 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 203..205 (old row 203 and everything below shifts
# down by 3, so old 203-212 become 206-215).
$ws.Range("A203:A205").EntireRow.Insert()

# Shared values for all three new rows.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$fecha = 44931
$codreg = 5
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria = "Damasco"
$variedad = "Patterson"
$unidad = "$/caja 15 kilos"
$origen = "Provincia de San Felipe de Aconcagua"
$kgUnidad = 15

# Row 203: Calidad Especial
$ws.Range("A203").Value = $mercadoId
$ws.Range("B203").Value = $mercado
$ws.Range("C203").Value = $region
$ws.Range("D203").Value = $fecha
$ws.Range("E203").Value = $codreg
$ws.Range("F203").Value = $tipo
$ws.Range("G203").Value = $productoId
$ws.Range("H203").Value = $producto
$ws.Range("I203").Value = $categoriaId
$ws.Range("J203").Value = $categoria
$ws.Range("K203").Value = $variedad
$ws.Range("L203").Value = "Especial"
$ws.Range("M203").Value = 56
$ws.Range("N203").Value = 14000
$ws.Range("O203").Value = 14000
$ws.Range("P203").Value = 14000
$ws.Range("Q203").Value = $unidad
$ws.Range("R203").Value = $origen
$ws.Range("S203").Value = 933
$ws.Range("T203").Value = $kgUnidad

# Row 204: Calidad Primera
$ws.Range("A204").Value = $mercadoId
$ws.Range("B204").Value = $mercado
$ws.Range("C204").Value = $region
$ws.Range("D204").Value = $fecha
$ws.Range("E204").Value = $codreg
$ws.Range("F204").Value = $tipo
$ws.Range("G204").Value = $productoId
$ws.Range("H204").Value = $producto
$ws.Range("I204").Value = $categoriaId
$ws.Range("J204").Value = $categoria
$ws.Range("K204").Value = $variedad
$ws.Range("L204").Value = "Primera"
$ws.Range("M204").Value = 68
$ws.Range("N204").Value = 12000
$ws.Range("O204").Value = 12000
$ws.Range("P204").Value = 12000
$ws.Range("Q204").Value = $unidad
$ws.Range("R204").Value = $origen
$ws.Range("S204").Value = 800
$ws.Range("T204").Value = $kgUnidad

# Row 205: Calidad Segunda
$ws.Range("A205").Value = $mercadoId
$ws.Range("B205").Value = $mercado
$ws.Range("C205").Value = $region
$ws.Range("D205").Value = $fecha
$ws.Range("E205").Value = $codreg
$ws.Range("F205").Value = $tipo
$ws.Range("G205").Value = $productoId
$ws.Range("H205").Value = $producto
$ws.Range("I205").Value = $categoriaId
$ws.Range("J205").Value = $categoria
$ws.Range("K205").Value = $variedad
$ws.Range("L205").Value = "Segunda"
$ws.Range("M205").Value = 60
$ws.Range("N205").Value = 10000
$ws.Range("O205").Value = 10000
$ws.Range("P205").Value = 10000
$ws.Range("Q205").Value = $unidad
$ws.Range("R205").Value = $origen
$ws.Range("S205").Value = 667
$ws.Range("T205").Value = $kgUnidad
